$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 76668
$ws.Range("C3").Value = 109655

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1%"
$ws.Range("E3").Style = "Normal"

$ws.Range("C4").Value = 257687

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2%"
$ws.Range("E4").Style = "Normal"

$ws.Range("C6").Value = 39800
